$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 and add rows 3-10 with new order data (store, order#, date, time, shipping method, note)

$ws.Range("A2").Value = 'موبایل جوان هشتگرد'
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '1785'
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '1402-10-22'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '18:02'
$ws.Range("E2").Value = 'پست'
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = 'فاکتور'

$ws.Range("A3").Value = 'موبایل جوان هشتگرد'
$ws.Range("B3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '1402-10-22'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '18:02'
$ws.Range("E3").Value = 'پست'
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = 'تست'

$ws.Range("A4").Value = 'موبایل جوان هشتگرد'
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '121212'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '1402-10-22'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '18:02'
$ws.Range("E4").Value = 'پست'
$ws.Range("F4").NumberFormat = "@"

$ws.Range("A5").Value = 'موبوپلاس ارومیه'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '89'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '1402-10-22'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '12:27'
$ws.Range("E5").Value = 'پست'
$ws.Range("F5").NumberFormat = "@"

$ws.Range("A6").Value = 'فروشگاه طنین موزیک'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '54545'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '1402-10-22'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '17:58'
$ws.Range("E6").Value = 'پست'
$ws.Range("F6").NumberFormat = "@"

$ws.Range("A7").Value = 'فروشگاه کامپیوتر R+'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '4565'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = '1402-10-19'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '17:55'
$ws.Range("E7").Value = 'پست'
$ws.Range("F7").NumberFormat = "@"

$ws.Range("A8").Value = 'های استور'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '12'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '1402-10-01'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '17:34'
$ws.Range("E8").Value = 'پست'
$ws.Range("F8").NumberFormat = "@"

$ws.Range("A9").Value = 'مصطفی ساری'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '1'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '1402-10-22'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '17:18'
$ws.Range("E9").Value = 'پست'
$ws.Range("F9").NumberFormat = "@"

$ws.Range("A10").Value = 'موبایل جوان هشتگرد'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '11111'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = '1402-10-22'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17:15'
$ws.Range("E10").Value = 'پست'
$ws.Range("F10").NumberFormat = "@"
